$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = -5.869799999999997
$ws.Range("C9").Value = -11.90290000000001
$ws.Range("D11").Value = -8.438099999999999
$ws.Range("C18").Value = -14.16259999999998
$ws.Range("C20").Value = -13.41669999999998
